$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the 25/04/2020 case count (row 52)
$ws.Range("B52").Value = 1835

# Insert a new row at 54 for 26/04/2020, shifting subsequent rows down
$ws.Rows.Item(54).Insert()
$ws.Range("A54").Value = "26/04/2020"
$ws.Range("B54").Value = 186
